# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" for the two entries that were just
# handed off again (ace6ff0c-... and the dependent c6f6881d-... entry which
# shares its handoff timestamp), on both the zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D11").Value = "2016-03-09 09:42:05"
$zhcn.Range("D12").Value = "2016-03-09 09:42:05"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D11").Value = "2016-03-09 09:42:08"
$dede.Range("D12").Value = "2016-03-09 09:42:08"
